$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "First name"
$ws.Range("B1").Value = "Last name"
$ws.Range("C1").Value = "Initials"
$ws.Range("D1").Value = "Pers.No."
$ws.Range("E1").Value = "Job Family"
$ws.Range("F1").Value = "Employee Subgroup"
$ws.Range("G1").Value = "Job"
$ws.Range("H1").Value = "Cost Center"
$ws.Range("I1").Value = "Init.Entry"
$ws.Range("J1").Value = "Pers.no. Superior"
$ws.Range("K1").Value = "Pers.no. Mentor"
$ws.Range("L1").Value = "Date of Birth"
$ws.Range("M1").Value = "Personnel Subarea"

# --- Column A / B (mentee, mentor, bestmentor / developer) ---
$ws.Range("A2").Value = "Mentee"
$ws.Range("B2").Value = "Developer"
$ws.Range("A3").Value = "Mentor"
$ws.Range("B3").Value = "Developer"
$ws.Range("A4").Value = "BestMentor"
$ws.Range("B4").Value = "Developer"

# --- Column C (initials) ---
$ws.Range("C2").Value = "TZJE"
$ws.Range("C3").Value = "DZLI"
$ws.Range("C4").Value = "SNKK"

# --- Column D (Pers.No.) ---
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3

# --- Column E (Job Family) ---
$ws.Range("E2").Value = "Project Development"
$ws.Range("E3").Value = "Project Development"
$ws.Range("E4").Value = "Project Development"

# --- Column F (Employee Subgroup) ---
$ws.Range("F2").Value = "Staff"
$ws.Range("F3").Value = "Staff"
$ws.Range("F4").Value = "Staff"

# --- Column G (Job) ---
$ws.Range("G2").Value = "L4 (Professional)"
$ws.Range("G3").Value = "L3 (Expert)"
$ws.Range("G4").Value = "L7 (Leader)"

# --- Column H (Cost Center) ---
$ws.Range("H2").Value = "Delivery"
$ws.Range("H3").Value = "Delivery"
$ws.Range("H4").Value = "Delivery"

# --- Column I (Init.Entry date) ---
$ws.Range("I2").Value = [DateTime]"2014-05-05"
$ws.Range("I3").Value = [DateTime]"2014-05-05"
$ws.Range("I4").Value = [DateTime]"2014-05-05"

# --- Column J (Pers.no. Superior) ---
# J2 already holds the text "601505" from the original workbook; leave it
# untouched so it keeps its Text type instead of being re-parsed as a Number.
$ws.Range("J3").Value = 123
$ws.Range("J4").Value = 456

# --- Column K (Pers.no. Mentor) ---
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 456
$ws.Range("K4").Value = 123

# --- Column L (Date of Birth) ---
$ws.Range("L2").Value = [DateTime]"1983-09-23"
$ws.Range("L3").Value = [DateTime]"1991-07-24"
$ws.Range("L4").Value = [DateTime]"1980-07-24"

# --- Column M (Personnel Subarea) ---
$ws.Range("M2").Value = "Lodz"
$ws.Range("M3").Value = "Lodz"
$ws.Range("M4").Value = "Lodz"

# --- Selection change ---
$ws.Range("D11").Select()
